# edit.ps1 - applies the commit's changes:
#  1. Remove <w:lastRenderedPageBreak/> from the "Narzedzie do ekstrakcji..." heading run.
#  2. Remove <w:lastRenderedPageBreak/> from the "Pytania" heading run.
#  3. After the "Jak dziala regex..." paragraph, add a new "Zrodla" Heading-1 paragraph
#     (carrying the lastRenderedPageBreak that moved down) plus a "Context - <url>" paragraph.

$d = $word.ActiveDocument

function Remove-LastRenderedPageBreak($searchText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $para = $r.Paragraphs(1)
    $pr = $para.Range
    $owxml = $pr.WordOpenXML
    if ($owxml -match '(<w:p\b.*?</w:p>)') {
        $pFull = $matches[1]
        $newXml = $pFull -replace '<w:lastRenderedPageBreak\s*/>', ''
        $pr.InsertXML($newXml) | Out-Null
    }
}

# 1 & 2: strip the two lastRenderedPageBreak markers.
Remove-LastRenderedPageBreak("Narzędzie do ekstrakcji tekstu z danych:")
Remove-LastRenderedPageBreak("Pytania")

# 3: append the new "Zrodla" section at the end of the document body.
$end = $d.Content
$end.Collapse(0)

$newSectionXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Nagwek1"/><w:rPr><w:lang w:val="pl-PL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pl-PL"/></w:rPr><w:lastRenderedPageBreak/><w:t>Źródła</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Context -</w:t></w:r><w:r><w:t>https://www.sohamkamani.com/golang/exec-shell-command/</w:t></w:r></w:p>'

$end.InsertXML($newSectionXml) | Out-Null

Write-Output "Edit applied."
